$d = $word.ActiveDocument

# 1. Update the letter date
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address line into a street line and a
#    city/state/zip line, followed by a new blank paragraph.
$d.Content.Find.Execute("4177 Stewart LN, Santa Clara, CA 95054", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4177 Stewart LN^pSanta Clara,, CA 95054^p", 2)

# 3. Remove the two now-superfluous blank paragraphs that used to sit
#    right after "Board of Directors" (find it dynamically since the
#    address split shifted every subsequent paragraph index by two).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Board of Directors*") {
        $targetIndex = $i
        break
    }
}
$d.Paragraphs($targetIndex + 1).Range.Delete()
$d.Paragraphs($targetIndex + 1).Range.Delete()
